$wb = $excel.ActiveWorkbook

# There will now be two vaccination dictionaries, one long and one short.
# Rename the existing "Vaccination" sheets to "Vaccination_long" (the
# "_xlnm._FilterDatabase" defined name that points at the Vaccination sheet
# is tracked by name and updates automatically when the sheet is renamed).
$wb.Worksheets.Item("Vaccination").Name = "Vaccination_long"
$wb.Worksheets.Item("Vaccination_options").Name = "Vaccination_long_options"
$wb.Worksheets.Item("Vaccination_settings").Name = "Vaccination_long_settings"

# Remembered selection on the (renamed) Vaccination_long_settings sheet.
$wb.Worksheets.Item("Vaccination_long_settings").Range("P47").Select()

# Switch the active tab back to "Mortality" (first sheet) with its
# remembered selection, instead of "Vaccination_long".
$wsMort = $wb.Worksheets.Item("Mortality")
$wsMort.Activate()
$wsMort.Range("D164").Select()
